# December 28 task check-in:
# Insert a new "Resource" column (with value "Test" on both data rows)
# immediately before the existing URL column (old column C), shifting
# everything from column C onward one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 0) Stash the hyperlink cells' original formatting (bold/underline
#    "Inconsolata" style) in a scratch area so it can be restored after
#    Hyperlinks.Add() below re-applies Excel's default Hyperlink style.
$ws.Range("C3").Copy()
$ws.Range("A100").PasteSpecial(-4122)

# 1) Insert a new blank column at C (old C "URL" -> D, etc.)
$ws.Columns("C:C").Insert()

# 2) Populate the new column with header + values
$ws.Range("C1").Value = "Resource"
$ws.Range("C2").Value = "Test"
$ws.Range("C3").Value = "Test"

# 3) The conditional formatting that used to live on G2 needs to track the
#    shifted cell, now H2.
$fc = $ws.Range("G2").FormatConditions.Item(1)
$fc.Formula1 = "=LEN(TRIM(H2))>0"
$fc.ModifyAppliesToRange($ws.Range("H2"))

# 4) The two hyperlinks used to anchor on column C (now D) - rebuild them
#    on the correct shifted cells, preserving their original targets/text.
#    D4 has no actual cell value (it never did), so add the hyperlink and
#    then drop the placeholder cell/row the Add() call creates for it.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D3"), "https://live.virtualandemo.com/api/pets/findByTags?tags=grey")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://live.virtualandemo.com/api/pets/findByTags?tags=grey", "", "", "https://live.virtualandemo.com/api/pets/findByTags?tags=grey")

# 4b) Hyperlinks.Add() stamps the default blue/underline "Hyperlink" style
#     on D3; put the original custom formatting back, then clear the
#     scratch cell used to stash it (before the row shift below).
$ws.Range("A100").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("A100").Clear()

# 4c) Drop the placeholder row the second Add() created for D4.
$ws.Rows("4:4").Delete()

# 5) Restore the current selection to match the edited workbook.
$ws.Range("C3").Select()
